# Apply weekly crime-data refresh to the 70th Precinct CompStat workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text runs)
# ---------------------------------------------------------------------------
# "Volume 32   Number  17" -> "... Number  18"
$ws.Range("A8").Characters(21, 2).Text = "18"

# "Report Covering the Week  4/21/2025  Through  4/27/2025"
#   -> "...4/28/2025  Through  5/4/2025"
$ws.Range("C9").Characters(27, 9).Text = "4/28/2025"
$ws.Range("C9").Characters(47, 9).Text = "5/4/2025"

# ---------------------------------------------------------------------------
# Column H width bumped to match column E (7.433768)
# ---------------------------------------------------------------------------
$ws.Columns("H").ColumnWidth = $ws.Columns("E").ColumnWidth

# ---------------------------------------------------------------------------
# Helper: cells that must become the text placeholder "0" (shared string,
# style 13) - copy format+value from a donor cell that already has that
# exact representation (e.g. C14).
# ---------------------------------------------------------------------------
$placeholderDonor = $ws.Range("C14")

function Set-Placeholder($addr) {
    $placeholderDonor.Copy($ws.Range($addr))
}

# Helper: cells that must switch FROM the text placeholder TO a plain
# integer value - set the value then copy the number format from a donor
# numeric cell in the same style family.
function Set-IntFromPlaceholder($addr, $value, $donorAddr) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.NumberFormat = $ws.Range($donorAddr).NumberFormat
}

# ---------------------------------------------------------------------------
# Row 15 - Robbery
# ---------------------------------------------------------------------------
Set-IntFromPlaceholder "C15" 2 "C16"
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 57.142857142857
$ws.Range("L15").Value = 83.333333333333
$ws.Range("M15").Value = 57.142857142857
$ws.Range("N15").Value = -74.418604651162

# ---------------------------------------------------------------------------
# Row 16 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -13.333333333333
$ws.Range("I16").Value = 66
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = 11.864406779661
$ws.Range("L16").Value = 26.923076923076
$ws.Range("M16").Value = -48.031496062992
$ws.Range("N16").Value = -90.476190476190

# ---------------------------------------------------------------------------
# Row 17 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 145
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = 106.896551724138
$ws.Range("L17").Value = 52.542372881355
$ws.Range("M17").Value = 48.760330578512
$ws.Range("N17").Value = -36.395759717314

# ---------------------------------------------------------------------------
# Row 18 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-Placeholder "C18"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 71.428571428571
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = 13.953488372093
$ws.Range("L18").Value = -14.035087719298
$ws.Range("M18").Value = -46.153846153846
$ws.Range("N18").Value = -95.205479452054

# ---------------------------------------------------------------------------
# Row 19 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -46.666666666666
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -40.350877192982
$ws.Range("I19").Value = 148
$ws.Range("J19").Value = 184
$ws.Range("K19").Value = -19.565217391304
$ws.Range("L19").Value = -24.489795918367
$ws.Range("M19").Value = -25.628140703517
$ws.Range("N19").Value = -60.321715817694

# ---------------------------------------------------------------------------
# Row 20 - TOTAL (non-bold)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 4
Set-IntFromPlaceholder "D20" 8 "C20"
Set-IntFromPlaceholder "E20" -50 "E18"
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -7.692307692307
$ws.Range("I20").Value = 41
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = -8.888888888888
$ws.Range("L20").Value = 20.588235294117
$ws.Range("M20").Value = -42.253521126760
$ws.Range("N20").Value = -95.170789163722

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 115
$ws.Range("H21").Value = 7.826086956521
$ws.Range("I21").Value = 495
$ws.Range("J21").Value = 425
$ws.Range("K21").Value = 16.470588235294
$ws.Range("L21").Value = 6.451612903225
$ws.Range("M21").Value = -20.032310177706
$ws.Range("N21").Value = -84.903934126258

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-Placeholder "F22"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -75

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 14.814814814814
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = -20.300751879699
$ws.Range("I24").Value = 504
$ws.Range("J24").Value = 638
$ws.Range("K24").Value = -21.003134796238
$ws.Range("L24").Value = -15.436241610738
$ws.Range("M24").Value = 30.232558139534

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 43
$ws.Range("H25").Value = -53.260869565217
$ws.Range("I25").Value = 263
$ws.Range("J25").Value = 397
$ws.Range("K25").Value = -33.753148614609
$ws.Range("L25").Value = -24.857142857142

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 27.272727272727
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 28.571428571428
$ws.Range("I26").Value = 223
$ws.Range("J26").Value = 214
$ws.Range("K26").Value = 4.205607476635
$ws.Range("L26").Value = 10.945273631840
$ws.Range("M26").Value = -8.979591836734

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 16.666666666666
$ws.Range("L27").Value = 7.692307692307

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -28
$ws.Range("L28").Value = -5.263157894736

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-Placeholder "C29"
$ws.Range("N29").Value = -92.857142857142

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-Placeholder "C30"
$ws.Range("N30").Value = -90.909090909090

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 0
